# Update the "票房/人气" (column F) figures on the "展览" and "全部类型"
# sheets to reflect newly generated output (gh-pages output @ 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 367
$ws1.Range("F3").Value  = 2276
$ws1.Range("F6").Value  = 5159
$ws1.Range("F10").Value = 235
$ws1.Range("F11").Value = 41
$ws1.Range("F12").Value = 220
$ws1.Range("F16").Value = 4227
$ws1.Range("F17").Value = 761
$ws1.Range("F18").Value = 769
$ws1.Range("F25").Value = 109
$ws1.Range("F26").Value = 599
$ws1.Range("F28").Value = 36
$ws1.Range("F29").Value = 1071
$ws1.Range("F31").Value = 2698
$ws1.Range("F33").Value = 82

# --- Sheet "全部类型" -------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 367
$ws4.Range("F3").Value  = 2276
$ws4.Range("F6").Value  = 5159
$ws4.Range("F10").Value = 235
$ws4.Range("F11").Value = 41
$ws4.Range("F12").Value = 220
$ws4.Range("F16").Value = 4227
$ws4.Range("F17").Value = 761
$ws4.Range("F18").Value = 769
$ws4.Range("F25").Value = 109
$ws4.Range("F26").Value = 599
$ws4.Range("F29").Value = 36
$ws4.Range("F30").Value = 1071
$ws4.Range("F32").Value = 2698
$ws4.Range("F34").Value = 82

$wb.Save()
